$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab: "linear_VTRF" -> "slinear_VTRF"
$ws.Name = "slinear_VTRF"

# C77:C84 were entered as whole-number percentages (e.g. 10.87) instead
# of the fractional values (0.1087) used everywhere else in column C.
$ws.Range("C77").Value = 0.1087
$ws.Range("C78").Value = 0.0826
$ws.Range("C79").Value = 0.109
$ws.Range("C80").Value = 0.0646
$ws.Range("C81").Value = 0.1945
$ws.Range("C82").Value = 0.1362
$ws.Range("C83").Value = 0.0848
$ws.Range("C84").Value = 0.1353

# Scroll the sheet down and move the active selection to C85, matching
# the saved view state.
$ws.Activate()
$ws.Range("C85").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
